$wb = $excel.ActiveWorkbook

$wsStarting = $wb.Worksheets.Item("Starting")
$wsAttack   = $wb.Worksheets.Item("Attack")
$wsSkill    = $wb.Worksheets.Item("Skill")
$wsPower    = $wb.Worksheets.Item("Power")

# Set the new card-name cells in the exact order the new shared strings were
# introduced, so the shared string table is built up in the same sequence.
$wsAttack.Range("A18").Value   = "Star Surge"
$wsSkill.Range("A9").Value     = "Rewind"
$wsPower.Range("A5").Value     = "Memento"
$wsSkill.Range("A14").Value    = "Parallel Universe"
$wsSkill.Range("A27").Value    = "Wormhole"
$wsSkill.Range("A8").Value     = "Time Dilation"
$wsPower.Range("A13").Value    = "Temporal Paradox"
$wsAttack.Range("A19").Value   = "Ring Singularity"
$wsAttack.Range("A22").Value   = "Cosmic Binding"
$wsSkill.Range("A29").Value    = "Tempered Fate"
$wsAttack.Range("A24").Value   = "Essence Flux"
$wsAttack.Range("A3").Value    = "Mystic Shot"
$wsAttack.Range("A26").Value   = "Arcane Barrage"
$wsSkill.Range("A18").Value    = "Mystic Shift"
$wsSkill.Range("A26").Value    = "Mimic"
$wsPower.Range("A8").Value     = "Distorted Reality"
$wsPower.Range("A7").Value     = "Inherit Wisdom"
$wsSkill.Range("A20").Value    = "Parallel Convergence"
$wsAttack.Range("A27").Value   = "Time Bomb"
$wsStarting.Range("A5").Value  = "Time Warp"
$wsStarting.Range("A4").Value  = "Mystic Blast"

# --- Selections / active cells per sheet ---
$wsAttack.Range("A8").Select()
$wsSkill.Range("A9").Select()
$wsPower.Range("A12").Select()

# Starting becomes the active (visible) tab/sheet, with A5 selected there last
$wsStarting.Activate()
$wsStarting.Range("A5").Select()
